$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.132146835327148
$ws.Range("B1").Value = 3.021248340606689
$ws.Range("C1").Value = 2.350072860717773
$ws.Range("D1").Value = 2.243273496627808
$ws.Range("E1").Value = 2.198022127151489
